$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.827.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.99%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.654.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.25%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.69%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  +1.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.668.04"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.27%  "

$ws.Range("E11").Value = "  +0.99%  "

$ws.Range("E12").Value = "  +0.97%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.135"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.71%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.122.90"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.729.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.92%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.653.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.15%  "

$ws.Range("E18").Value = "  +1.52%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "345.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.96%  "

$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.77%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.417"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.52%  "

$ws.Range("E26").Value = "  -1.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.88%  "

$ws.Range("E29").Value = "  +2.73%  "

$ws.Range("E30").Value = "  -0.04%  "

$ws.Range("E31").Value = "  +1.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.82%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.06"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.56%  "

$ws.Range("E36").Value = "  +2.31%  "

$ws.Range("E37").Value = "  +0.93%  "

$ws.Range("E38").Value = "  -0.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.828"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "294.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.30%  "

$ws.Range("E41").Value = "  +2.13%  "

$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.607"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0543"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.05%  "

$ws.Range("E46").Value = "  -0.35%  "

$ws.Range("E47").Value = "  -1.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.983.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.17%  "

$ws.Range("E49").Value = "  +1.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.82%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.93%  "
